# Horarios actualizados Linea 141 - scrape refresh 04:10:26 -> 04:50:21
$wb = $excel.ActiveWorkbook

$newStamp = "04:50:21"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newStamp"
$ws1.Range("A3").Value = "Total filas: 12"

# The earliest departure (row 6, 81_EL PELIGRO) has already happened by the
# new scrape time, so it drops off the top of the list and everything else
# shifts up one row.
$ws1.Rows.Item(6).Delete()

# Full refreshed data set (rows 6-17), new scrape timestamp + recomputed
# "Minutos" countdown; five new arrivals appended at the bottom.
$sheet1Data = @(
    @("04:53", "11_ETCHEVERRY",    3),
    @("05:17", "17_ROMERO",       27),
    @("05:22", "23_HERNANDEZ",    32),
    @("05:44", "14_ABASTO",       54),
    @("05:47", "17_ROMERO",       57),
    @("06:01", "16_SANTA ANA",    71),
    @("06:09", "10_OLMOS",        79),
    @("06:16", "215A_EL PATO",    86),
    @("06:30", "23_HERNANDEZ",   100),
    @("06:34", "11_ETCHEVERRY",  104),
    @("06:39", "17X38_ROMERO",   109),
    @("06:41", "16_SANTA ANA",   111)
)

$r = 6
foreach ($row in $sheet1Data) {
    $ws1.Cells.Item($r, 1).Value = $newStamp
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = "LP1912"
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newStamp"
$ws2.Range("A3").Value = "Total filas: 1"

$ws2.Range("A5").Value = "Hora_Scrap"
$ws2.Range("B5").Value = "Hora_Llegada"
$ws2.Range("C5").Value = "Linea"
$ws2.Range("D5").Value = "Minutos"
$ws2.Range("E5").Value = "Parada"

$ws2.Range("A6").Value = $newStamp
$ws2.Range("B6").Value = "06:16"
$ws2.Range("C6").Value = "215A_EL PATO"
$ws2.Range("D6").Value = 86
$ws2.Range("E6").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newStamp"
